$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '90.684.99'
$ws.Range("E2").Value = '  +1.60%  '

# Row 3
$ws.Range("D3").Value = '3.153.71'
$ws.Range("E3").Value = '  +4.02%  '

# Row 4
$ws.Range("E4").Value = '  -0.42%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.65'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.82%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '626.06'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.38%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.14'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +29.51%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.369'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.46%  '

# Row 9
$ws.Range("E9").Value = '  -0.17%  '

# Row 10
$ws.Range("D10").Value = '3.151.63'
$ws.Range("E10").Value = '  +4.06%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.763'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +15.30%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.201'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +8.07%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000244'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.60%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.66'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +5.56%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '34.73'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +9.44%  '

# Row 16
$ws.Range("D16").Value = '90.328.64'
$ws.Range("E16").Value = '  +2.24%  '

# Row 17
$ws.Range("D17").Value = '3.730.37'
$ws.Range("E17").Value = '  +3.60%  '

# Row 18
$ws.Range("D18").Value = '3.180.22'
$ws.Range("E18").Value = '  +4.87%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.66'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +9.79%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.29'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +7.17%  '

# Row 21
$ws.Range("B21").Value = 'BitcoinCash'
$ws.Range("C21").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '463.00'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +9.27%  '

# Row 22
$ws.Range("B22").Value = 'PEPE'
$ws.Range("C22").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0000210'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.58%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.08'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +12.25%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.29'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +6.25%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.84'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +8.91%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '89.17'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +7.24%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.11'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +4.03%  '

# Row 28
$ws.Range("D28").Value = '3.313.22'
$ws.Range("E28").Value = '  +3.60%  '

# Row 29
$ws.Range("E29").Value = '  +0.00%  '

# Row 30
$ws.Range("B30").Value = 'Binance-PegBSC-USD'
$ws.Range("C30").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.03'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.64%  '

# Row 31
$ws.Range("B31").Value = 'Cronos'
$ws.Range("C31").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.162'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.08%  '

# Row 32
$ws.Range("B32").Value = 'InternetComputer(DFINITY)'
$ws.Range("C32").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '9.22'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +13.25%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '27.20'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +19.73%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '518.42'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.45%  '

# Row 35
$ws.Range("B35").Value = 'Stellar'
$ws.Range("C35").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.183'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +34.16%  '

# Row 36
$ws.Range("B36").Value = 'dogwifhat'
$ws.Range("C36").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.62'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.14%  '

# Row 37
$ws.Range("B37").Value = 'PancakeSwap'
$ws.Range("C37").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.91'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +7.10%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.142'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +9.88%  '

# Row 39
$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.86'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +4.07%  '

# Row 40
$ws.Range("E40").Value = '  +5.07%  '

# Row 41
$ws.Range("B41").Value = 'Hedera'
$ws.Range("C41").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0875'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +29.28%  '

# Row 42
$ws.Range("B42").Value = 'WhiteBITCoin'
$ws.Range("C42").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '22.20'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.13%  '

# Row 43
$ws.Range("E43").Value = '  -0.13%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.411'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +14.42%  '

# Row 45
$ws.Range("B45").Value = 'Stacks'
$ws.Range("C45").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.95'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +7.38%  '

# Row 46
$ws.Range("B46").Value = 'USDe'
$ws.Range("C46").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.00'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.04%  '

# Row 47
$ws.Range("B47").Value = 'Filecoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.61'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +13.95%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '45.83'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +5.63%  '

# Row 49
$ws.Range("B49").Value = 'Monero'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '148.16'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.63%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.34'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +11.33%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.667'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +14.84%  '

